$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 249.54546
$ws.Range("I4").Value = 138.33333
$ws.Range("K4").Value = 138.33333
$ws.Range("M4").Value = -24.33332999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 506.125
$ws.Range("I12").Value = 435.57144
$ws.Range("K12").Value = 435.57144
$ws.Range("M12").Value = -265.57144

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3590.5334
$ws.Range("I62").Value = 1654.3334
$ws.Range("K62").Value = 1654.3334
$ws.Range("M62").Value = -1030.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3590.5334
$ws.Range("I65").Value = 1654.3334
$ws.Range("K65").Value = 8271.666999999999
$ws.Range("M65").Value = -5151.666999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 68944.125
$ws.Range("J133").Value = 68944.125
$ws.Range("L133").Value = 68944.125
$ws.Range("N133").Value = -79064.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2857.5334
$ws.Range("I135").Value = 2804.1428
$ws.Range("J135").Value = 2904.25
$ws.Range("K135").Value = 25237.2852
$ws.Range("L135").Value = 26138.25
$ws.Range("M135").Value = -22702.2852
$ws.Range("N135").Value = -31208.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1673432.4
$ws.Range("I8").Value = 2503898.5
$ws.Range("K8").Value = 2503898.5
$ws.Range("M8").Value = -2503754.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 5393.6665
$ws.Range("I22").Value = 5393.6665
$ws.Range("K22").Value = 5393.6665
$ws.Range("M22").Value = -5094.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2787.5583
$ws.Range("I32").Value = 2310.3013
$ws.Range("J32").Value = 11497.5
$ws.Range("K32").Value = 2310.3013
$ws.Range("L32").Value = 11497.5
$ws.Range("M32").Value = -2023.3013
$ws.Range("N32").Value = -12071.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13892768
$ws.Range("I74").Value = 22224838
$ws.Range("J74").Value = 5986
$ws.Range("K74").Value = 22224838
$ws.Range("L74").Value = 5986
$ws.Range("M74").Value = -22223964
$ws.Range("N74").Value = -7734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 13892768
$ws.Range("I77").Value = 22224838
$ws.Range("J77").Value = 5986
$ws.Range("K77").Value = 111124190
$ws.Range("L77").Value = 29930
$ws.Range("M77").Value = -111119822
$ws.Range("N77").Value = -38666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5545.7837
$ws.Range("I132").Value = 5343.857
$ws.Range("J132").Value = 6174
$ws.Range("K132").Value = 16031.571
$ws.Range("L132").Value = 18522
$ws.Range("M132").Value = -13501.571
$ws.Range("N132").Value = -23582

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 10026
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 10026
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 10026
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -10530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2810.158
$ws.Range("I134").Value = 2810.158
$ws.Range("K134").Value = 8430.474
$ws.Range("M134").Value = -5895.474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 28328.5
$ws.Range("J43").Value = 28328.5
$ws.Range("L43").Value = 28328.5
$ws.Range("N43").Value = -28696.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 28328.5
$ws.Range("J101").Value = 28328.5
$ws.Range("L101").Value = 28328.5
$ws.Range("N101").Value = -34818.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2914.0356
$ws.Range("I132").Value = 1889.6842
$ws.Range("J132").Value = 5076.5557
$ws.Range("K132").Value = 5669.0526
$ws.Range("L132").Value = 15229.6671
$ws.Range("M132").Value = -3139.0526
$ws.Range("N132").Value = -20289.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4041.8235
$ws.Range("I134").Value = 2407
$ws.Range("K134").Value = 7221
$ws.Range("M134").Value = -4686

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2606236
$ws.Range("I107").Value = 2511.5
$ws.Range("J107").Value = 5209960.5
$ws.Range("K107").Value = 7534.5
$ws.Range("L107").Value = 15629881.5
$ws.Range("M107").Value = -5614.5
$ws.Range("N107").Value = -15633721.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1211.5454
$ws.Range("J113").Value = 1289.3125
$ws.Range("L113").Value = 3867.9375
$ws.Range("N113").Value = -8207.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3031.6667
$ws.Range("I132").Value = 1023.7143
$ws.Range("J132").Value = 4035.6428
$ws.Range("K132").Value = 9213.4287
$ws.Range("L132").Value = 36320.7852
$ws.Range("M132").Value = -6683.4287
$ws.Range("N132").Value = -41380.7852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4500002
$ws.Range("I14").Value = 4
$ws.Range("J14").Value = 9000000
$ws.Range("K14").Value = 4
$ws.Range("L14").Value = 9000000
$ws.Range("M14").Value = 164
$ws.Range("N14").Value = -9000336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 48265.383
$ws.Range("I20").Value = 34502.5
$ws.Range("K20").Value = 34502.5
$ws.Range("M20").Value = -34257.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 55038.555
$ws.Range("I24").Value = 54005.668
$ws.Range("K24").Value = 54005.668
$ws.Range("M24").Value = -53832.668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 50000
$ws.Range("K60").Value = 50000
$ws.Range("M60").Value = -49425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 112110010
$ws.Range("J17").Value = 112110010
$ws.Range("L17").Value = 112110010
$ws.Range("N17").Value = -112110350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2385.6924
$ws.Range("I30").Value = 1699.2
$ws.Range("J30").Value = 2814.75
$ws.Range("K30").Value = 1699.2
$ws.Range("L30").Value = 2814.75
$ws.Range("M30").Value = -1591.2
$ws.Range("N30").Value = -3030.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 2686.2
$ws.Range("I35").Value = 857.75
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 857.75
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = -521.75
$ws.Range("N35").Value = -10672

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7913.7144
$ws.Range("J132").Value = 9786
$ws.Range("L132").Value = 29358
$ws.Range("N132").Value = -34418

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10373.685
$ws.Range("I136").Value = 4156.5557
$ws.Range("J136").Value = 15969.1
$ws.Range("K136").Value = 12469.6671
$ws.Range("L136").Value = 47907.3
$ws.Range("M136").Value = -9919.667099999999
$ws.Range("N136").Value = -53007.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 15000
$ws.Range("J32").Value = 15000
$ws.Range("L32").Value = 15000
$ws.Range("N32").Value = -15634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2430.25
$ws.Range("I122").Value = 1813.3928
$ws.Range("K122").Value = 5440.178400000001
$ws.Range("M122").Value = -2990.178400000001
